# Auto-update draw results: append the latest Pick 3 draw as a new row
# right after the last populated row of the results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1048576").End(-4162).Row   # xlUp
$newRow = $lastRow + 1

$rng = $ws.Range("A" + $newRow + ":E" + $newRow)
$rng.NumberFormat = "@"   # keep date/phase strings as text, not auto-converted

$ws.Cells.Item($newRow, 1).Value = "2025-12-10"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251210"
$ws.Cells.Item($newRow, 4).Value = "5-2-5"
$ws.Cells.Item($newRow, 5).Value = "2025-12-10T21:45:10.263+04:00"

$rng.Style = "Normal"   # drop the temporary text format so it matches the rest of the sheet
